$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 923.37177
$ws.Range("J17").Value = 865.44446
$ws.Range("L17").Value = 2596.33338
$ws.Range("N17").Value = -2932.33338
$ws.Range("H19").Value = 4387218.5
$ws.Range("I19").Value = 8772530
$ws.Range("J19").Value = 1907.3334
$ws.Range("K19").Value = 8772530
$ws.Range("L19").Value = 1907.3334
$ws.Range("M19").Value = -8772355
$ws.Range("N19").Value = -2257.3334
$ws.Range("H115").Value = 1030.2727
$ws.Range("I115").Value = 833.3
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 2499.9
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -932.8999999999996
$ws.Range("N115").Value = -12134
$ws.Range("H116").Value = 536189.8
$ws.Range("I116").Value = 1669800.1
$ws.Range("K116").Value = 1669800.1
$ws.Range("M116").Value = -1666358.1
$ws.Range("H132").Value = 32264352
$ws.Range("I132").Value = 37043030
$ws.Range("J132").Value = 8251.5
$ws.Range("K132").Value = 111129090
$ws.Range("L132").Value = 24754.5
$ws.Range("M132").Value = -111126560
$ws.Range("N132").Value = -29814.5
$ws.Range("H135").Value = 552.4737
$ws.Range("I135").Value = 577.2
$ws.Range("J135").Value = 525
$ws.Range("K135").Value = 5194.8
$ws.Range("L135").Value = 4725
$ws.Range("M135").Value = -2659.8
$ws.Range("N135").Value = -9795
$ws.Range("H138").Value = 5409.88
$ws.Range("I138").Value = 1155.9375
$ws.Range("J138").Value = 6220.155
$ws.Range("K138").Value = 3467.8125
$ws.Range("L138").Value = 18660.465
$ws.Range("M138").Value = 1672.1875
$ws.Range("N138").Value = -28940.465
$ws.Range("H141").Value = 27833.395
$ws.Range("I141").Value = 30616.588
$ws.Range("J141").Value = 4176.25
$ws.Range("K141").Value = 91849.764
$ws.Range("L141").Value = 12528.75
$ws.Range("M141").Value = -86669.764
$ws.Range("N141").Value = -22888.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3807.1738
$ws.Range("I32").Value = 3444.5166
$ws.Range("J32").Value = 6224.8887
$ws.Range("K32").Value = 3444.5166
$ws.Range("L32").Value = 6224.8887
$ws.Range("M32").Value = -3157.5166
$ws.Range("N32").Value = -6798.8887
$ws.Range("H45").Value = 1172.4
$ws.Range("I45").Value = 1006
$ws.Range("K45").Value = 1006
$ws.Range("M45").Value = -629
$ws.Range("H74").Value = 3238.975
$ws.Range("I74").Value = 3264.875
$ws.Range("J74").Value = 3135.375
$ws.Range("K74").Value = 3264.875
$ws.Range("L74").Value = 3135.375
$ws.Range("M74").Value = -2390.875
$ws.Range("N74").Value = -4883.375
$ws.Range("H77").Value = 3238.975
$ws.Range("I77").Value = 3264.875
$ws.Range("J77").Value = 3135.375
$ws.Range("K77").Value = 16324.375
$ws.Range("L77").Value = 15676.875
$ws.Range("M77").Value = -11956.375
$ws.Range("N77").Value = -24412.875
$ws.Range("H132").Value = 2145.827
$ws.Range("I132").Value = 1278.2703
$ws.Range("J132").Value = 4285.8
$ws.Range("K132").Value = 3834.810899999999
$ws.Range("L132").Value = 12857.4
$ws.Range("M132").Value = -1304.810899999999
$ws.Range("N132").Value = -17917.4
$ws.Range("H137").Value = 36442.5
$ws.Range("J137").Value = 39791.43
$ws.Range("L137").Value = 39791.43
$ws.Range("N137").Value = -49991.43
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1697.221
$ws.Range("I134").Value = 1059.0476
$ws.Range("J134").Value = 3445.261
$ws.Range("K134").Value = 3177.142800000001
$ws.Range("L134").Value = 10335.783
$ws.Range("M134").Value = -642.1428000000005
$ws.Range("N134").Value = -15405.783
$ws.Range("H137").Value = 44461
$ws.Range("J137").Value = 49326.25
$ws.Range("L137").Value = 49326.25
$ws.Range("N137").Value = -59526.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15390767
$ws.Range("I99").Value = 40004380
$ws.Range("J99").Value = 7258.75
$ws.Range("K99").Value = 40004380
$ws.Range("L99").Value = 7258.75
$ws.Range("M99").Value = -40002882
$ws.Range("N99").Value = -10254.75
$ws.Range("H112").Value = 27425
$ws.Range("J112").Value = 27425
$ws.Range("L112").Value = 27425
$ws.Range("N112").Value = -30379
$ws.Range("H126").Value = 15390767
$ws.Range("I126").Value = 40004380
$ws.Range("J126").Value = 7258.75
$ws.Range("K126").Value = 120013140
$ws.Range("L126").Value = 21776.25
$ws.Range("M126").Value = -120010670
$ws.Range("N126").Value = -26716.25
$ws.Range("H137").Value = 49975.293
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 49975.293
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 49975.293
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -60175.293
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 256
$ws.Range("J75").Value = 256
$ws.Range("L75").Value = 768
$ws.Range("N75").Value = -2764
$ws.Range("H78").Value = 256
$ws.Range("J78").Value = 256
$ws.Range("L78").Value = 2304
$ws.Range("N78").Value = -12288
$ws.Range("H97").Value = 320.55554
$ws.Range("I97").Value = 252.5
$ws.Range("J97").Value = 456.66666
$ws.Range("K97").Value = 757.5
$ws.Range("L97").Value = 1369.99998
$ws.Range("M97").Value = -261.5
$ws.Range("N97").Value = -2361.99998
$ws.Range("H113").Value = 523.1515000000001
$ws.Range("I113").Value = 450.79413
$ws.Range("K113").Value = 1352.38239
$ws.Range("M113").Value = 817.61761
$ws.Range("H131").Value = 824.28986
$ws.Range("I131").Value = 483.33334
$ws.Range("J131").Value = 875.43335
$ws.Range("K131").Value = 1450.00002
$ws.Range("L131").Value = 2626.30005
$ws.Range("M131").Value = 3589.99998
$ws.Range("N131").Value = -12706.30005
$ws.Range("H132").Value = 2556.8333
$ws.Range("I132").Value = 1091.5
$ws.Range("J132").Value = 2849.9
$ws.Range("K132").Value = 9823.5
$ws.Range("L132").Value = 25649.1
$ws.Range("M132").Value = -7293.5
$ws.Range("N132").Value = -30709.1
$ws.Range("H136").Value = 3460.111
$ws.Range("J136").Value = 4758.25
$ws.Range("L136").Value = 14274.75
$ws.Range("N136").Value = -24474.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 32228
$ws.Range("J46").Value = 32228
$ws.Range("L46").Value = 32228
$ws.Range("N46").Value = -32540
$ws.Range("H102").Value = 2421.6956
$ws.Range("I102").Value = 1861.4
$ws.Range("J102").Value = 3472.25
$ws.Range("K102").Value = 1861.4
$ws.Range("L102").Value = 3472.25
$ws.Range("M102").Value = -239.4000000000001
$ws.Range("N102").Value = -6716.25
$ws.Range("H122").Value = 3740.25
$ws.Range("I122").Value = 1662.8
$ws.Range("K122").Value = 4988.4
$ws.Range("M122").Value = -2538.4
$ws.Range("H126").Value = 2876.0203
$ws.Range("I126").Value = 2895.0205
$ws.Range("J126").Value = 1014
$ws.Range("K126").Value = 8685.0615
$ws.Range("L126").Value = 3042
$ws.Range("M126").Value = -6215.0615
$ws.Range("N126").Value = -7982
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4782.2354
$ws.Range("J7").Value = 6362.5
$ws.Range("L7").Value = 6362.5
$ws.Range("N7").Value = -6586.5
$ws.Range("H61").Value = 1098.2963
$ws.Range("I61").Value = 993.9091
$ws.Range("J61").Value = 1557.6
$ws.Range("K61").Value = 993.9091
$ws.Range("L61").Value = 1557.6
$ws.Range("M61").Value = -791.9091
$ws.Range("N61").Value = -1961.6
$ws.Range("H93").Value = 4446065
$ws.Range("J93").Value = 1958.2307
$ws.Range("L93").Value = 1958.2307
$ws.Range("N93").Value = -4454.2307
$ws.Range("H113").Value = 1098.2963
$ws.Range("I113").Value = 993.9091
$ws.Range("J113").Value = 1557.6
$ws.Range("K113").Value = 993.9091
$ws.Range("L113").Value = 1557.6
$ws.Range("M113").Value = 1176.0909
$ws.Range("N113").Value = -5897.6
$ws.Range("H122").Value = 6260.3335
$ws.Range("I122").Value = 2975
$ws.Range("J122").Value = 7455
$ws.Range("K122").Value = 8925
$ws.Range("L122").Value = 22365
$ws.Range("M122").Value = -6475
$ws.Range("N122").Value = -27265
$ws.Range("H126").Value = 4782.2354
$ws.Range("J126").Value = 6362.5
$ws.Range("L126").Value = 19087.5
$ws.Range("N126").Value = -24027.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 43200
$ws.Range("J112").Value = 43200
$ws.Range("L112").Value = 43200
$ws.Range("N112").Value = -46154
$ws.Range("H122").Value = 2735.6667
$ws.Range("I122").Value = 1669.2106
$ws.Range("J122").Value = 5268.5
$ws.Range("K122").Value = 5007.6318
$ws.Range("L122").Value = 15805.5
$ws.Range("M122").Value = -2557.6318
$ws.Range("N122").Value = -20705.5
$ws.Range("H132").Value = 5377640.5
$ws.Range("I132").Value = 527.1395
$ws.Range("J132").Value = 17546896
$ws.Range("K132").Value = 1581.4185
$ws.Range("L132").Value = 52640688
$ws.Range("M132").Value = 948.5815
$ws.Range("N132").Value = -52645748
$ws.Range("H136").Value = 1214.8281
$ws.Range("I136").Value = 702.6667
$ws.Range("J136").Value = 2192.5908
$ws.Range("K136").Value = 2108.0001
$ws.Range("L136").Value = 6577.7724
$ws.Range("M136").Value = 441.9998999999998
$ws.Range("N136").Value = -11677.7724
